# Applies the Halicarnassus_Profits market-data refresh produced by the scheduled runner.
# For each changed leve row, currentAveragePrice/NQ/HQ (H,I,J), LevePrice NQ/HQ (K,L),
# and LeveProfit NQ/HQ (M,N) are updated to the freshly-pulled market values. Cells that
# no longer have a computed profit (e.g. price data dropped to 0) are cleared entirely,
# matching upstream, rather than left stale.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 212.26666
$ws.Range("I9").Value = 275.8
$ws.Range("J9").Value = 85.2
$ws.Range("K9").Value = 275.8
$ws.Range("L9").Value = 85.2
$ws.Range("M9").Value = -106.8
$ws.Range("N9").Value = -423.2
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
# Row 55
$ws.Range("H55").Value = 276.85715
$ws.Range("I55").Value = 188.22223
$ws.Range("J55").Value = 436.4
$ws.Range("K55").Value = 188.22223
$ws.Range("L55").Value = 436.4
$ws.Range("M55").Value = 25.77777
$ws.Range("N55").Value = -864.4
# Row 113
$ws.Range("H113").Value = 4743.143
$ws.Range("I113").Value = 2991
$ws.Range("J113").Value = 5035.1665
$ws.Range("K113").Value = 2991
$ws.Range("L113").Value = 5035.1665
$ws.Range("M113").Value = 263
$ws.Range("N113").Value = -11543.1665
# Row 137
$ws.Range("H137").Value = 3207.5
$ws.Range("I137").Value = 1708.6666
$ws.Range("J137").Value = 4492.2144
$ws.Range("K137").Value = 5125.9998
$ws.Range("L137").Value = 13476.6432
$ws.Range("M137").Value = -2575.9998
$ws.Range("N137").Value = -18576.6432

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 950.5
$ws.Range("I2").Value = 913.86664
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 913.86664
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -800.86664
$ws.Range("N2").Value = -1726
# Row 45
$ws.Range("H45").Value = 2611.25
$ws.Range("I45").Value = 2243.375
$ws.Range("K45").Value = 2243.375
$ws.Range("M45").Value = -1866.375
# Row 105
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
# Row 110
$ws.Range("H110").Value = 788.26666
$ws.Range("I110").Value = 786.53845
$ws.Range("J110").Value = 799.5
$ws.Range("K110").Value = 786.53845
$ws.Range("L110").Value = 799.5
$ws.Range("M110").Value = 1258.46155
$ws.Range("N110").Value = -4889.5
# Row 116
$ws.Range("H116").Value = 950.5
$ws.Range("I116").Value = 913.86664
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 913.86664
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1380.13336
$ws.Range("N116").Value = -6088
# Row 122
$ws.Range("H122").Value = 3639.8
$ws.Range("I122").Value = 3424.75
$ws.Range("K122").Value = 10274.25
$ws.Range("M122").Value = -7824.25

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 950.5
$ws.Range("I3").Value = 913.86664
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 913.86664
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -799.86664
$ws.Range("N3").Value = -1728
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
# Row 39
$ws.Range("H39").Value = 16633.334
$ws.Range("J39").Value = 16633.334
$ws.Range("L39").Value = 16633.334
$ws.Range("N39").Value = -17411.334
# Row 95
$ws.Range("H95").Value = 27436.666
$ws.Range("J95").Value = 27436.666
$ws.Range("L95").Value = 27436.666
$ws.Range("N95").Value = -32928.666

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2171.6316
$ws.Range("J58").Value = 2360.9
$ws.Range("L58").Value = 2360.9
$ws.Range("N58").Value = -2766.9
# Row 132
$ws.Range("H132").Value = 4585.7334
$ws.Range("I132").Value = 4378.6
$ws.Range("K132").Value = 13135.8
$ws.Range("M132").Value = -10605.8
# Row 134
$ws.Range("H134").Value = 999.5
$ws.Range("I134").Value = 999
$ws.Range("K134").Value = 2997
$ws.Range("M134").Value = -462
# Row 136
$ws.Range("H136").Value = 2171.6316
$ws.Range("J136").Value = 2360.9
$ws.Range("L136").Value = 7082.700000000001
$ws.Range("N136").Value = -12182.7

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2749.3333
$ws.Range("I68").Value = 2748.5
$ws.Range("K68").Value = 8245.5
$ws.Range("M68").Value = -7434.5
# Row 71
$ws.Range("H71").Value = 2749.3333
$ws.Range("I71").Value = 2748.5
$ws.Range("K71").Value = 24736.5
$ws.Range("M71").Value = -20680.5
# Row 103
$ws.Range("H103").Value = 614.1667
$ws.Range("I103").Value = 383.75
$ws.Range("J103").Value = 1075
$ws.Range("K103").Value = 1151.25
$ws.Range("L103").Value = 3225
$ws.Range("M103").Value = -272.25
$ws.Range("N103").Value = -4983
# Row 131
$ws.Range("H131").Value = 1616
$ws.Range("I131").Value = 1013
$ws.Range("K131").Value = 3039
$ws.Range("M131").Value = 2001

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 5875.75
$ws.Range("J14").Value = 12252.5
$ws.Range("L14").Value = 12252.5
$ws.Range("N14").Value = -12588.5
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
# Row 102
$ws.Range("H102").Value = 1836.1428
$ws.Range("I102").Value = 1719.6364
$ws.Range("J102").Value = 2263.3333
$ws.Range("K102").Value = 1719.6364
$ws.Range("L102").Value = 2263.3333
$ws.Range("M102").Value = -97.63640000000009
$ws.Range("N102").Value = -5507.3333
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
# Row 126
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 7729.9
$ws.Range("I100").Value = 4324.75
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 4324.75
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -3783.75
$ws.Range("N100").Value = -11082
# Row 122
$ws.Range("H122").Value = 3633.3333
$ws.Range("I122").Value = 3450
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10350
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -7900
$ws.Range("N122").Value = -16900
# Row 132
$ws.Range("H132").Value = 3499.3333
$ws.Range("J132").Value = 3499.5
$ws.Range("L132").Value = 10498.5
$ws.Range("N132").Value = -15558.5
# Row 135
$ws.Range("H135").Value = 100429
$ws.Range("J135").Value = 100429
$ws.Range("L135").Value = 100429
$ws.Range("N135").Value = -110569

$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 87604.5
$ws.Range("J80").Value = 87604.5
$ws.Range("L80").Value = 87604.5
$ws.Range("N80").Value = -89600.5
# Row 83
$ws.Range("H83").Value = 87604.5
$ws.Range("J83").Value = 87604.5
$ws.Range("L83").Value = 262813.5
$ws.Range("N83").Value = -272797.5
# Row 107
$ws.Range("H107").Value = 2464.1428
$ws.Range("I107").Value = 2464.1428
$ws.Range("K107").Value = 7392.428400000001
$ws.Range("M107").Value = -5472.428400000001
